$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows for the new "Blockee" Action Figure products,
# right before the existing "Charizard + Charmander Duo Keepley Figure" row (old row 119).
$ws.Rows("119:122").Insert()

# Column A (Product names) first so the new shared strings are appended
# in product-name order, matching how Excel would record them.
$ws.Range("A119").Value = "Charizard Blockee Figure"
$ws.Range("A120").Value = "Greninja Blockee Figure"
$ws.Range("A121").Value = "Meowscrada Blockee Figure"
$ws.Range("A122").Value = "Ceruledge Blockee Figure"

# Column D (Image file names) next.
$ws.Range("D119").Value = "charizard blockee.jpg"
$ws.Range("D120").Value = "greninja blockee.jpg"
$ws.Range("D121").Value = "grass blockee.jpg"
$ws.Range("D122").Value = "sword blockee.jpg"

# Column B (Price).
$ws.Range("B119").Value = 3750
$ws.Range("B120").Value = 3200
$ws.Range("B121").Value = 3200
$ws.Range("B122").Value = 3350

# Column C (Category) - reuses the existing "Action Figure" shared string.
$ws.Range("C119").Value = "Action Figure"
$ws.Range("C120").Value = "Action Figure"
$ws.Range("C121").Value = "Action Figure"
$ws.Range("C122").Value = "Action Figure"

# Update the selection to match the post-edit cursor position.
$ws.Range("D122").Select()
